$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$app = $excel

# --- Append two new survey rows (39, 40) dated 2023-10-26 / serial 45225 -------------
#
# The simplest way to get pixel/style-perfect new rows is to duplicate an existing
# row's formatting (Insert a copy of row 38, which already carries the date-number
# format in column A and the percentage-with-"Porcentaje"-style formats used across
# columns C/D/F/H and C/E/G/I/J) and then overwrite just the values.

# Row 39: clone row 38's formatting
$ws.Rows("38:38").Copy()
$ws.Rows("39:39").Insert(-4121)   # xlShiftDown

# Row 40: clone row 38's formatting again
$ws.Rows("38:38").Copy()
$ws.Rows("40:40").Insert(-4121)   # xlShiftDown

$app.CutCopyMode = $false

# Column B ("Encuestadora") needs a different fill per pollster; reuse the exact
# style already present on existing "Guarumo" / "Mosqueteros" rows instead of
# re-deriving it, so no new cell styles are introduced.
$ws.Range("B5").Copy($ws.Range("B39"))   # B5 is an existing "Guarumo" row
$ws.Range("B3").Copy($ws.Range("B40"))   # B3 is an existing "Mosqueteros" row
$app.CutCopyMode = $false

# --- Row 39 values: Guarumo, 2023-10-26 ---
$ws.Range("A39").Value2 = 45225
$ws.Range("C39").Value2 = 0.322
$ws.Range("D39").Value2 = 0.331
$ws.Range("E39").Value2 = 0.097
$ws.Range("F39").Value2 = 0.065
$ws.Range("G39").Value2 = 0.048
$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
# "Otros" is the left-over share; compute it the same way the rest of the sheet
# was built (1 - sum of the other shares) so the stored value matches bit-for-bit.
$ws.Range("J39").Formula = "=1-SUM(C39:I39)"

# --- Row 40 values: Mosqueteros, 2023-10-26 ---
$ws.Range("A40").Value2 = 45225
$ws.Range("C40").Value2 = 0.3272
$ws.Range("D40").Value2 = 0.1915
$ws.Range("E40").Value2 = 0.0115
$ws.Range("F40").Value2 = 0.2634
$ws.Range("G40").Value2 = 0.1673
$ws.Range("H40").Value2 = 0.0033
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Formula = "=1-SUM(C40:I40)"

# Freeze the computed "Otros" formulas down to static values, matching every
# other cell in the sheet (no live formulas outside the external-link rows).
$rng = $ws.Range("J39:J40")
$rng.Copy()
$rng.PasteSpecial(-4163)   # xlPasteValues
$app.CutCopyMode = $false

# Match the saved selection state from the edit.
[void]$ws.Range("G47").Select()

Write-Host "Added rows 39-40 (Guarumo, Mosqueteros) for 2023-10-26"
